# "update to figure 2"
# 1) Bump the auto-generated datetimeFigureOut footer date (3/26/2019 -> 3/27/2019)
#    everywhere it appears: the slide master and every slide layout.
# 2) Rework Figure 2 (slide 5): widen/shift the "Matching Variables (demographics)"
#    textbox and relabel it "(Demographics)"; relabel the unit-of-observation
#    textbox from "(Missions)" to "(Nonprofits)".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {}
        if ($isDate -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "3/26/2019") {
                $tr.Text = "3/27/2019"
            }
        }
    }
}

# -- Slide master --
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# -- Every slide layout under the master --
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}

# -- Figure 2 (slide 5): the two callout textboxes live inside "Group 8" --
$slide5 = $p.Slides.Item(5)
$group = $slide5.Shapes.Item(1)

$tbDemographics = $null
$tbUnitOfObs = $null
for ($i = 1; $i -le $group.GroupItems.Count; $i++) {
    $item = $group.GroupItems.Item($i)
    if ($item.Name -eq "TextBox 4") { $tbDemographics = $item }
    if ($item.Name -eq "TextBox 5") { $tbUnitOfObs = $item }
}

# Reposition / resize the "Matching Variables (demographics)" textbox:
#   off  x: 8648487 -> 8638068 EMU
#   ext cx: 1648272 -> 1669111 EMU
# (EMU / 12700 = points; nudge by half an EMU so the runtime rounds to the
#  exact target EMU value instead of truncating.)
$tbDemographics.Left = (8638068.5 / 12700)
$tbDemographics.Width = (1669111.5 / 12700)

# "(demographics)" -> "(Demographics)" split into two runs: "(Demographics" + ")"
$para = $tbDemographics.TextFrame.TextRange.Paragraphs(2)
$firstPart = $para.Characters(1, 13)
$firstPart.Text = "(Demographics"

# "Unit of Observation (Missions)" -> "Unit of Observation " + "(Nonprofits)"
$para2 = $tbUnitOfObs.TextFrame.TextRange.Paragraphs(1)
$secondPart = $para2.Characters(21, 10)
$secondPart.Text = "(Nonprofits)"

Write-Output "done"
